$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12, pushing the existing rows 12-26 down to
# 13-27 (this also grows the sheet dimension from A1:R26 to A1:R27 and
# carries the date number-format down into the new D12 cell).
$ws.Rows(12).Insert()

# Populate the newly inserted row 12 with the new weekly record.
$ws.Range("A12").Value2 = 7
$ws.Range('B12').Value2 = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range('C12').Value2 = 'Ñuble'
$ws.Range("D12").Value2 = 44161
$ws.Range("E12").Value2 = 16
$ws.Range("F12").Value2 = 100112022
$ws.Range('G12').Value2 = 'Arveja Verde'
$ws.Range('H12').Value2 = 'Sin especificar'
$ws.Range('I12').Value2 = 'Primera'
$ws.Range("J12").Value2 = 34
$ws.Range("K12").Value2 = 19500
$ws.Range("L12").Value2 = 20000
$ws.Range("M12").Value2 = 19735
$ws.Range('N12').Value2 = '$/saco 25 kilos'
$ws.Range('O12').Value2 = 'Región del Maule'
$ws.Range("P12").Value2 = 789
$ws.Range("Q12").Value2 = 25
$ws.Range('R12').Value2 = 'Hortaliza'
